# Updates cryptos list values per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'34.358.88"
$ws.Cells.Item(2, 5).Value = "'  -0.86%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.799.69"
$ws.Cells.Item(3, 5).Value = "'  -1.04%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "'  -0.05%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'226.34"
$ws.Cells.Item(5, 5).Value = "'  -0.87%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.596"
$ws.Cells.Item(6, 5).Value = "'  +2.30%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "'  +0.00%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'36.20"
$ws.Cells.Item(8, 5).Value = "'  +4.47%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.294"
$ws.Cells.Item(9, 5).Value = "'  -2.38%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.0688"
$ws.Cells.Item(10, 5).Value = "'  -1.82%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0962"
$ws.Cells.Item(11, 5).Value = "'  +0.95%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'2.056.70"
$ws.Cells.Item(12, 5).Value = "'  -1.22%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'11.30"
$ws.Cells.Item(13, 5).Value = "'  -0.85%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'1.794.35"
$ws.Cells.Item(14, 5).Value = "'  -1.49%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.637"
$ws.Cells.Item(15, 5).Value = "'  -1.03%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'34.347.39"
$ws.Cells.Item(16, 5).Value = "'  -0.88%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'4.41"
$ws.Cells.Item(17, 5).Value = "'  +1.76%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'69.58"
$ws.Cells.Item(18, 5).Value = "'  +0.57%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'243.08"
$ws.Cells.Item(19, 5).Value = "'  -1.26%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'0.0₃0783"
$ws.Cells.Item(20, 5).Value = "'  -2.38%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'11.40"
$ws.Cells.Item(21, 5).Value = "'  -1.31%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "'  +0.07%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'4.15"
$ws.Cells.Item(23, 5).Value = "'  -0.71%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'2.25"
$ws.Cells.Item(24, 5).Value = "'  +7.18%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'170.60"
$ws.Cells.Item(25, 5).Value = "'  -1.77%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'8.01"
$ws.Cells.Item(26, 5).Value = "'  +7.02%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'17.31"
$ws.Cells.Item(27, 5).Value = "'  +3.04%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "'  +0.99%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "'  -0.02%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "'Filecoin"
$ws.Cells.Item(30, 3).Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).Value = "'3.81"
$ws.Cells.Item(30, 5).Value = "'  -0.86%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'3.93"
$ws.Cells.Item(31, 5).Value = "'  -1.50%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "'PancakeSwap"
$ws.Cells.Item(32, 3).Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(32, 4).Value = "'1.23"
$ws.Cells.Item(32, 5).Value = "'  -0.88%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "'  -2.51%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.79"
$ws.Cells.Item(34, 5).Value = "'  -2.91%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'1.373.40"
$ws.Cells.Item(35, 5).Value = "'  -2.38%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.654"
$ws.Cells.Item(36, 5).Value = "'  -3.98%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "'  -1.70%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'2.36"
$ws.Cells.Item(38, 5).Value = "'  -10.26%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "'  -2.72%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'81.94"
$ws.Cells.Item(40, 5).Value = "'  -2.61%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "'  +0.79%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "'  -3.05%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'0.943"
$ws.Cells.Item(43, 5).Value = "'  -0.76%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "'  +5.34%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'13.44"
$ws.Cells.Item(45, 5).Value = "'  -2.01%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "'  -4.28%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'1.960.06"
$ws.Cells.Item(47, 5).Value = "'  -1.18%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'5.85"
$ws.Cells.Item(48, 5).Value = "'  -3.65%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "'  +0.12%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'102.35"
$ws.Cells.Item(50, 5).Value = "'  -2.85%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "'BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51, 4).Value = "'0.0₆0123"
$ws.Cells.Item(51, 5).Value = "'  -5.07%  "
